# Gym log Q1 2024 - append new workout entries (rows 715-746)
# and set the active selection to match the author's final cursor position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 715
$ws.Cells.Item(715, 1).Value = 45461
$ws.Cells.Item(715, 2).Value = "Cardio"
$ws.Cells.Item(715, 3).Value = "Run"
$ws.Cells.Item(715, 4).Value = 2
$ws.Cells.Item(715, 9).Value = 0.25
$ws.Cells.Item(715, 10).Value = "Out"

# Row 716
$ws.Cells.Item(716, 1).Value = 45461
$ws.Cells.Item(716, 2).Value = "Cardio"

# Row 717
$ws.Cells.Item(717, 1).Value = 45462
$ws.Cells.Item(717, 2).Value = "Cardio"
$ws.Cells.Item(717, 3).Value = "Run"
$ws.Cells.Item(717, 4).Value = 1
$ws.Cells.Item(717, 9).Value = 0.38
$ws.Cells.Item(717, 10).Value = "Gym"

# Row 718
$ws.Cells.Item(718, 1).Value = 45462
$ws.Cells.Item(718, 2).Value = "Biceps"
$ws.Cells.Item(718, 3).Value = "Bicep curl"
$ws.Cells.Item(718, 4).Value = 1
$ws.Cells.Item(718, 5).Value = 10
$ws.Cells.Item(718, 6).Value = 13

# Row 719
$ws.Cells.Item(719, 1).Value = 45462
$ws.Cells.Item(719, 2).Value = "Biceps"
$ws.Cells.Item(719, 3).Value = "Bicep curl"
$ws.Cells.Item(719, 4).Value = 1
$ws.Cells.Item(719, 5).Value = 6
$ws.Cells.Item(719, 6).Value = 15

# Row 720
$ws.Cells.Item(720, 1).Value = 45462
$ws.Cells.Item(720, 2).Value = "Triceps"
$ws.Cells.Item(720, 3).Value = "Kneeling dip"
$ws.Cells.Item(720, 4).Value = 4
$ws.Cells.Item(720, 5).Value = 6
$ws.Cells.Item(720, 6).Value = 0
$ws.Cells.Item(720, 7).Value = "BW"

# Row 721
$ws.Cells.Item(721, 1).Value = 45462
$ws.Cells.Item(721, 2).Value = "Chest"
$ws.Cells.Item(721, 3).Value = "Incline chest press bench"
$ws.Cells.Item(721, 4).Value = 4
$ws.Cells.Item(721, 5).Value = 8
$ws.Cells.Item(721, 6).Value = 30

# Row 722
$ws.Cells.Item(722, 1).Value = 45463
$ws.Cells.Item(722, 2).Value = "Cardio"
$ws.Cells.Item(722, 3).Value = "Run"
$ws.Cells.Item(722, 4).Value = 1
$ws.Cells.Item(722, 9).Value = 0.5
$ws.Cells.Item(722, 10).Value = "Gym"

# Row 723
$ws.Cells.Item(723, 1).Value = 45463
$ws.Cells.Item(723, 2).Value = "Legs"
$ws.Cells.Item(723, 3).Value = "Leg extension"
$ws.Cells.Item(723, 4).Value = 1
$ws.Cells.Item(723, 5).Value = 10
$ws.Cells.Item(723, 6).Value = 52

# Row 724
$ws.Cells.Item(724, 1).Value = 45463
$ws.Cells.Item(724, 2).Value = "Legs"
$ws.Cells.Item(724, 3).Value = "Leg extension"
$ws.Cells.Item(724, 4).Value = 3
$ws.Cells.Item(724, 5).Value = 10
$ws.Cells.Item(724, 6).Value = 61

# Row 725
$ws.Cells.Item(725, 1).Value = 45463
$ws.Cells.Item(725, 2).Value = "Legs"
$ws.Cells.Item(725, 3).Value = "Prone leg curl"
$ws.Cells.Item(725, 4).Value = 4
$ws.Cells.Item(725, 5).Value = 8
$ws.Cells.Item(725, 6).Value = 54

# Row 726
$ws.Cells.Item(726, 1).Value = 45463
$ws.Cells.Item(726, 2).Value = "Legs"
$ws.Cells.Item(726, 3).Value = "Stair master"
$ws.Cells.Item(726, 4).Value = 4
$ws.Cells.Item(726, 6).Value = 0
$ws.Cells.Item(726, 7).Value = "BW"
$ws.Cells.Item(726, 8).Value = "1 set = 10 floors"

# Row 727
$ws.Cells.Item(727, 1).Value = 45467
$ws.Cells.Item(727, 2).Value = "Cardio"
$ws.Cells.Item(727, 3).Value = "Run"
$ws.Cells.Item(727, 4).Value = 0.5
$ws.Cells.Item(727, 9).Value = 0.5
$ws.Cells.Item(727, 10).Value = "Gym"

# Row 728
$ws.Cells.Item(728, 1).Value = 45467
$ws.Cells.Item(728, 2).Value = "Chest"
$ws.Cells.Item(728, 3).Value = "Incline press"
$ws.Cells.Item(728, 4).Value = 4
$ws.Cells.Item(728, 5).Value = 12
$ws.Cells.Item(728, 6).Value = 45

# Row 729
$ws.Cells.Item(729, 1).Value = 45467
$ws.Cells.Item(729, 2).Value = "Chest"
$ws.Cells.Item(729, 3).Value = "Pec fly"
$ws.Cells.Item(729, 4).Value = 4
$ws.Cells.Item(729, 5).Value = 11
$ws.Cells.Item(729, 6).Value = 66

# Row 730
$ws.Cells.Item(730, 1).Value = 45467
$ws.Cells.Item(730, 2).Value = "Triceps"
$ws.Cells.Item(730, 3).Value = "Kneeling dip"
$ws.Cells.Item(730, 4).Value = 1
$ws.Cells.Item(730, 5).Value = 5
$ws.Cells.Item(730, 6).Value = 0
$ws.Cells.Item(730, 7).Value = "BW"

# Row 731
$ws.Cells.Item(731, 1).Value = 45467
$ws.Cells.Item(731, 2).Value = "Triceps"
$ws.Cells.Item(731, 3).Value = "Kneeling dip"
$ws.Cells.Item(731, 4).Value = 1
$ws.Cells.Item(731, 5).Value = 5
$ws.Cells.Item(731, 6).Value = -20
$ws.Cells.Item(731, 7).Value = "BW"

# Row 732
$ws.Cells.Item(732, 1).Value = 45467
$ws.Cells.Item(732, 2).Value = "Triceps"
$ws.Cells.Item(732, 3).Value = "Kneeling dip"
$ws.Cells.Item(732, 4).Value = 2
$ws.Cells.Item(732, 5).Value = 6
$ws.Cells.Item(732, 6).Value = -27
$ws.Cells.Item(732, 7).Value = "BW"

# Row 733
$ws.Cells.Item(733, 1).Value = 45407
$ws.Cells.Item(733, 2).Value = "Cardio"
$ws.Cells.Item(733, 3).Value = "Run"
$ws.Cells.Item(733, 4).Value = 1
$ws.Cells.Item(733, 9).Value = 0.5
$ws.Cells.Item(733, 10).Value = "Gym"

# Row 734
$ws.Cells.Item(734, 1).Value = 45407
$ws.Cells.Item(734, 2).Value = "Back"
$ws.Cells.Item(734, 3).Value = "Low row"
$ws.Cells.Item(734, 4).Value = 1
$ws.Cells.Item(734, 5).Value = 12
$ws.Cells.Item(734, 6).Value = 39

# Row 735
$ws.Cells.Item(735, 1).Value = 45407
$ws.Cells.Item(735, 2).Value = "Back"
$ws.Cells.Item(735, 3).Value = "Low row"
$ws.Cells.Item(735, 4).Value = 1
$ws.Cells.Item(735, 5).Value = 12
$ws.Cells.Item(735, 6).Value = 45

# Row 736
$ws.Cells.Item(736, 1).Value = 45407
$ws.Cells.Item(736, 2).Value = "Back"
$ws.Cells.Item(736, 3).Value = "Low row"
$ws.Cells.Item(736, 4).Value = 2
$ws.Cells.Item(736, 5).Value = 8
$ws.Cells.Item(736, 6).Value = 52

# Row 737
$ws.Cells.Item(737, 1).Value = 45407
$ws.Cells.Item(737, 2).Value = "Triceps"
$ws.Cells.Item(737, 3).Value = "Seated dip"
$ws.Cells.Item(737, 4).Value = 4
$ws.Cells.Item(737, 5).Value = 14
$ws.Cells.Item(737, 6).Value = 75

# Row 738
$ws.Cells.Item(738, 1).Value = 45407
$ws.Cells.Item(738, 2).Value = "Triceps"
$ws.Cells.Item(738, 3).Value = "Triceps push down"
$ws.Cells.Item(738, 4).Value = 2
$ws.Cells.Item(738, 5).Value = 10
$ws.Cells.Item(738, 6).Value = 18

# Row 739
$ws.Cells.Item(739, 1).Value = 45470
$ws.Cells.Item(739, 2).Value = "Cardio"
$ws.Cells.Item(739, 3).Value = "Run"
$ws.Cells.Item(739, 4).Value = 0.5
$ws.Cells.Item(739, 9).Value = 0.75
$ws.Cells.Item(739, 10).Value = "Gym"

# Row 740
$ws.Cells.Item(740, 1).Value = 45470
$ws.Cells.Item(740, 2).Value = "Biceps"
$ws.Cells.Item(740, 3).Value = "Bicep Hammer curl"
$ws.Cells.Item(740, 4).Value = 1
$ws.Cells.Item(740, 5).Value = 10
$ws.Cells.Item(740, 6).Value = 13

# Row 741
$ws.Cells.Item(741, 1).Value = 45470
$ws.Cells.Item(741, 2).Value = "Biceps"
$ws.Cells.Item(741, 3).Value = "Bicep Hammer curl"
$ws.Cells.Item(741, 4).Value = 1
$ws.Cells.Item(741, 5).Value = 9
$ws.Cells.Item(741, 6).Value = 15

# Row 742
$ws.Cells.Item(742, 1).Value = 45470
$ws.Cells.Item(742, 2).Value = "Biceps"
$ws.Cells.Item(742, 3).Value = "Bicep Hammer curl"
$ws.Cells.Item(742, 4).Value = 2
$ws.Cells.Item(742, 5).Value = 6
$ws.Cells.Item(742, 6).Value = 15

# Row 743
$ws.Cells.Item(743, 1).Value = 45470
$ws.Cells.Item(743, 2).Value = "Triceps"
$ws.Cells.Item(743, 3).Value = "Kneeling dip"
$ws.Cells.Item(743, 4).Value = 4
$ws.Cells.Item(743, 5).Value = 8
$ws.Cells.Item(743, 6).Value = 0
$ws.Cells.Item(743, 7).Value = "BW"

# Row 744
$ws.Cells.Item(744, 1).Value = 45470
$ws.Cells.Item(744, 2).Value = "Chest"
$ws.Cells.Item(744, 3).Value = "Bench press"
$ws.Cells.Item(744, 4).Value = 2
$ws.Cells.Item(744, 5).Value = 8
$ws.Cells.Item(744, 6).Value = 50

# Row 745
$ws.Cells.Item(745, 1).Value = 45470
$ws.Cells.Item(745, 2).Value = "Chest"
$ws.Cells.Item(745, 3).Value = "Chest press"
$ws.Cells.Item(745, 4).Value = 2
$ws.Cells.Item(745, 5).Value = 8
$ws.Cells.Item(745, 6).Value = 55

# Row 746
$ws.Cells.Item(746, 1).Value = 45470
$ws.Cells.Item(746, 2).Value = "Cardio"
$ws.Cells.Item(746, 3).Value = "Virtual bike"
$ws.Cells.Item(746, 4).Value = 2

# Restore the selection/active cell as last left by the author
[void]$ws.Range("D740").Select()

Write-Host "Added rows 715-746; new dimension and shared strings updated."
